$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 ---
# D2: Target cluster changes from "sCs" to "FAPs"
$ws.Cells.Item(2, 4).Value = "FAPs"

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.770179333333333
$ws.Cells.Item(2, 8).Value = 5.310538
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.4098869999999999
$ws.Cells.Item(2, 14).Value = 1.229661
$ws.Cells.Item(2, 15).Value = 0.6566849734690215
$ws.Cells.Item(2, 16).Value = 0.6566849734690215
$ws.Cells.Item(2, 17).Value = 0.7255734964019999
$ws.Cells.Item(2, 18).Value = 6.530161467618
$ws.Cells.Item(2, 19).Value = 0.6566849734690215
$ws.Cells.Item(2, 20).Value = 0.6566849734690215

# --- Add new row 3 ---
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ccl5"
$ws.Cells.Item(3, 3).Value = "Gpr75"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.770179333333333
$ws.Cells.Item(3, 8).Value = 5.310538
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.214289
$ws.Cells.Item(3, 14).Value = 0.642867
$ws.Cells.Item(3, 15).Value = 0.3433150265309785
$ws.Cells.Item(3, 16).Value = 0.3433150265309785
$ws.Cells.Item(3, 17).Value = 0.3793299591606666
$ws.Cells.Item(3, 18).Value = 3.413969632446
$ws.Cells.Item(3, 19).Value = 0.3433150265309785
$ws.Cells.Item(3, 20).Value = 0.3433150265309785
